$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.389.18'
$ws.Range('E2').Value = '  -0.82%  '
$ws.Range('D3').Value = '1.709.33'
$ws.Range('E3').Value = '  -1.86%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9953'
$ws.Range('E4').Value = '  -0.38%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '238.86'
$ws.Range('E5').Value = '  -3.38%  '
$ws.Range('E6').Value = '  -0.37%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4894'
$ws.Range('E7').Value = '  -0.73%  '
$ws.Range('E8').Value = '  -3.99%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06162'
$ws.Range('E9').Value = '  -2.15%  '
$ws.Range('D10').Value = '1.715.69'
$ws.Range('E10').Value = '  -1.46%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06937'
$ws.Range('E11').Value = '  -1.88%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '15.47'
$ws.Range('E12').Value = '  -1.88%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.5971'
$ws.Range('E14').Value = '  -3.30%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '76.29'
$ws.Range('E15').Value = '  -2.62%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9960'
$ws.Range('E16').Value = '  -0.37%  '
$ws.Range('D17').Value = '26.261.55'
$ws.Range('E17').Value = '  -1.30%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9954'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000007056'
$ws.Range('E19').Value = '  -3.48%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '11.20'
$ws.Range('E20').Value = '  -3.25%  '
$ws.Range('D21').Value = '1.935.51'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.365'
$ws.Range('E22').Value = '  -5.08%  '
$ws.Range('E23').Value = '  -4.21%  '
$ws.Range('E24').Value = '  -5.01%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.74'
$ws.Range('E25').Value = '  -1.92%  '
$ws.Range('E26').Value = '  -2.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.386'
$ws.Range('E27').Value = '  -2.44%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.718'
$ws.Range('E28').Value = '  -2.75%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '105.02'
$ws.Range('E29').Value = '  -2.97%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '3.872'
$ws.Range('E30').Value = '  -4.56%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.07934'
$ws.Range('E31').Value = '  -1.64%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.589'
$ws.Range('E32').Value = '  -4.06%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04437'
$ws.Range('E33').Value = '  -4.08%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.589'
$ws.Range('E34').Value = '  -0.76%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9921'
$ws.Range('E35').Value = '  -2.87%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6104'
$ws.Range('E36').Value = '  -4.44%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.9513'
$ws.Range('E37').Value = '  +5.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.986'
$ws.Range('E38').Value = '  -4.34%  '
$ws.Range('E39').Value = '  -1.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.9950'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.01477'
$ws.Range('E41').Value = '  -1.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '99.27'
$ws.Range('E42').Value = '  -2.71%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '5.357'
$ws.Range('E43').Value = '  -1.41%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.3777'
$ws.Range('E44').Value = '  -4.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '6.784'
$ws.Range('E45').Value = '  -2.02%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1144'
$ws.Range('E46').Value = '  -3.64%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.05337'
$ws.Range('E47').Value = '  -1.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '30.30'
$ws.Range('E48').Value = '  -1.15%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '7.691'
$ws.Range('E49').Value = '  -1.86%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '50.83'
$ws.Range('E50').Value = '  -1.90%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9984'
$ws.Range('E51').Value = '  -0.45%  '
